# Updates the cryptos list figures (prices / 1h volume changes) and fixes
# the ordering/content of the Mantle / EnergySwap rows, per the scraped
# data refresh performed by the GitHub Actions job.

function Set-CellText($sheet, $ref, $text) {
    # Assign as text. A leading apostrophe forces Excel to treat values
    # that look numeric (e.g. "9.40", "54.80") as literal text so that
    # formatting (trailing zeros, grouping dots, etc.) is preserved
    # exactly as scraped. Resetting the style back to "Normal" afterwards
    # strips the quote-prefix flag that the apostrophe trick adds, so the
    # cell keeps the workbook's original (default) styling.
    $sheet.Range($ref).Value = "'" + $text
    $sheet.Range($ref).Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Subscript digits used inside a couple of very-small-price tickers
# (e.g. 0.0<sub>3</sub>0741). Built via [char] so the literal unicode
# glyphs don't need to be embedded directly in the script file.
$sub3 = [string][char]0x2083
$sub6 = [string][char]0x2086

Set-CellText $ws "D2" "26.682.40"
Set-CellText $ws "E2" "  +1.52%  "
Set-CellText $ws "D3" "1.638.63"
Set-CellText $ws "E3" "  +1.90%  "
Set-CellText $ws "E4" "  -0.24%  "
Set-CellText $ws "D5" "213.05"
Set-CellText $ws "E5" "  +0.08%  "
Set-CellText $ws "E6" "  -0.19%  "
Set-CellText $ws "E7" "  +1.24%  "
Set-CellText $ws "D9" "0.0624"
Set-CellText $ws "E9" "  +1.48%  "
Set-CellText $ws "E10" "  +4.77%  "
Set-CellText $ws "E11" "  +2.66%  "
Set-CellText $ws "D12" "1.864.21"
Set-CellText $ws "E12" "  +1.66%  "
Set-CellText $ws "D13" "1.642.83"
Set-CellText $ws "E13" "  +1.34%  "
Set-CellText $ws "D14" "4.05"
Set-CellText $ws "E14" "  +0.46%  "
Set-CellText $ws "E15" "  +1.83%  "
Set-CellText $ws "D16" "26.668.18"
Set-CellText $ws "E16" "  +1.38%  "
Set-CellText $ws "D17" "63.06"
Set-CellText $ws "E17" "  +1.42%  "
Set-CellText $ws "D18" ("0.0" + $sub3 + "0741")
Set-CellText $ws "E18" "  +1.71%  "
Set-CellText $ws "D19" "210.16"
Set-CellText $ws "E19" "  +3.67%  "
Set-CellText $ws "E20" "  -0.11%  "
Set-CellText $ws "E21" "  +0.63%  "
Set-CellText $ws "D22" "9.40"
Set-CellText $ws "E22" "  +0.81%  "
Set-CellText $ws "E23" "  +2.02%  "
Set-CellText $ws "D24" "1.93"
Set-CellText $ws "E24" "  +3.10%  "
Set-CellText $ws "D25" "145.95"
Set-CellText $ws "E25" "  +0.69%  "
Set-CellText $ws "E26" "  -0.24%  "
Set-CellText $ws "E27" "  -0.74%  "
Set-CellText $ws "D28" "6.72"
Set-CellText $ws "E28" "  +2.41%  "
Set-CellText $ws "E29" "  +1.27%  "
Set-CellText $ws "D30" "0.0518"
Set-CellText $ws "E31" "  -0.51%  "
Set-CellText $ws "E32" "  +0.56%  "
Set-CellText $ws "E33" "  +1.10%  "
Set-CellText $ws "E34" "  +0.49%  "
Set-CellText $ws "E35" "  -0.99%  "
Set-CellText $ws "D36" "1.170.82"
Set-CellText $ws "E36" "  +0.40%  "
Set-CellText $ws "E37" "  +0.88%  "
Set-CellText $ws "E38" "  +2.60%  "
Set-CellText $ws "E39" "  -0.13%  "
Set-CellText $ws "E40" "  -0.16%  "
Set-CellText $ws "E41" "  +0.72%  "
Set-CellText $ws "D42" "0.796"
Set-CellText $ws "E42" "  +1.76%  "
Set-CellText $ws "E43" "  +1.62%  "
Set-CellText $ws "D44" "1.773.47"
Set-CellText $ws "E44" "  +1.35%  "
Set-CellText $ws "D45" "92.33"
Set-CellText $ws "E45" "  +0.31%  "
Set-CellText $ws "E46" "  +2.20%  "
Set-CellText $ws "D47" ("0.0" + $sub6 + "0105")
Set-CellText $ws "E47" "  +5.26%  "
Set-CellText $ws "D48" "54.80"
Set-CellText $ws "E48" "  +0.90%  "

# Rows 50/51 swapped order: Mantle <-> EnergySwap.
Set-CellText $ws "B50" "EnergySwap"
Set-CellText $ws "C50" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-CellText $ws "D50" "7.57"
Set-CellText $ws "E50" "  +4.23%  "
Set-CellText $ws "B51" "Mantle"
Set-CellText $ws "C51" "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-CellText $ws "D51" "0.409"
Set-CellText $ws "E51" "  +0.34%  "
